# Apply the Fri Sep  6 07:34:37 UTC 2024 "cryptos" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.437.86"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "2.329.15"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'498.76"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'127.78"
$ws.Range("E6").Value = "  -4.10%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.536"
$ws.Range("E8").Value = "  -3.36%  "
$ws.Range("D9").Value = "2.328.68"
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("D10").Value = "'0.0974"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'4.81"
$ws.Range("E12").Value = "  +4.64%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "2.746.23"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "55.445.65"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "'21.48"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "2.397.41"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'9.85"
$ws.Range("E19").Value = "  -4.20%  "
$ws.Range("D20").Value = "'307.66"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'64.94"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("E27").Value = "  -4.37%  "
$ws.Range("D28").Value = "'7.10"
$ws.Range("E28").Value = "  -4.76%  "
$ws.Range("D29").Value = "'172.76"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").Value = "0.0₃0696"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'5.78"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("D36").Value = "'17.53"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "'3.61"
$ws.Range("E38").Value = "  -5.91%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").Value = "'36.09"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").Value = "'126.05"
$ws.Range("E43").Value = "  -5.60%  "
$ws.Range("D44").Value = "'4.76"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("D45").Value = "'0.555"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("D47").Value = "'234.86"
$ws.Range("E47").Value = "  -6.81%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("E49").Value = "  -3.87%  "
$ws.Range("D50").Value = "'16.50"
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").Value = "'0.952"
$ws.Range("E51").Value = "  +0.15%  "
